# Apply the "Updated cryptos list" refresh: new prices / 1h-volume deltas,
# plus two coin-identity swaps (rows 13/14 and 41/42 traded rank positions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "27.561.19"
$ws.Range('E2').Value = "  +5.50%  "

# Row 4
$ws.Range('E4').Value = "  +0.05%  "

# Row 5
$ws.Range('D5').Value = "'226.27"
$ws.Range('E5').Value = "  +3.55%  "

# Row 6
$ws.Range('D6').Value = "'0.5382"
$ws.Range('E6').Value = "  +2.72%  "

# Row 7
$ws.Range('E7').Value = "  +0.01%  "

# Row 8
$ws.Range('D8').Value = "'0.2675"
$ws.Range('E8').Value = "  +0.75%  "

# Row 9
$ws.Range('D9').Value = "'0.06613"
$ws.Range('E9').Value = "  +4.08%  "

# Row 10
$ws.Range('D10').Value = "'21.79"
$ws.Range('E10').Value = "  +6.23%  "

# Row 11
$ws.Range('D11').Value = "'0.07722"
$ws.Range('E11').Value = "  +0.30%  "

# Row 12
$ws.Range('D12').Value = "'4.614"
$ws.Range('E12').Value = "  -0.48%  "

# Row 13: Coin -> WrappedEther
$ws.Range('B13').Value = "WrappedEther"
$ws.Range('C13').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D13').Value = "1.725.51"
$ws.Range('E13').Value = "  +3.43%  "

# Row 14: Coin -> WrappedliquidstakedEther2.0
$ws.Range('B14').Value = "WrappedliquidstakedEther2.0"
$ws.Range('C14').Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D14').Value = "1.963.04"
$ws.Range('E14').Value = "  +4.29%  "

# Row 15
$ws.Range('D15').Value = "'0.5875"
$ws.Range('E15').Value = "  +4.64%  "

# Row 16
$ws.Range('D16').Value = "0.0₅8320"
$ws.Range('E16').Value = "  +1.73%  "

# Row 17
$ws.Range('D17').Value = "'68.05"
$ws.Range('E17').Value = "  +3.98%  "

# Row 18
$ws.Range('D18').Value = "27.581.73"
$ws.Range('E18').Value = "  +5.59%  "

# Row 19
$ws.Range('D19').Value = "'221.98"
$ws.Range('E19').Value = "  +15.37%  "

# Row 21
$ws.Range('D21').Value = "'4.741"
$ws.Range('E21').Value = "  +1.94%  "

# Row 22
$ws.Range('D22').Value = "'10.70"
$ws.Range('E22').Value = "  +1.93%  "

# Row 23
$ws.Range('D23').Value = "'6.105"
$ws.Range('E23').Value = "  +2.53%  "

# Row 24
$ws.Range('D24').Value = "'1.004"
$ws.Range('E24').Value = "  +0.05%  "

# Row 25
$ws.Range('D25').Value = "'148.21"
$ws.Range('E25').Value = "  +2.36%  "

# Row 26
$ws.Range('D26').Value = "'1.695"
$ws.Range('E26').Value = "  +12.03%  "

# Row 27
$ws.Range('D27').Value = "'0.1235"
$ws.Range('E27').Value = "  +3.23%  "

# Row 28
$ws.Range('D28').Value = "'7.402"
$ws.Range('E28').Value = "  +1.88%  "

# Row 29
$ws.Range('D29').Value = "'16.67"
$ws.Range('E29').Value = "  +4.40%  "

# Row 30
$ws.Range('D30').Value = "'0.05550"
$ws.Range('E30').Value = "  +1.83%  "

# Row 31
$ws.Range('D31').Value = "'1.305"
$ws.Range('E31').Value = "  +2.73%  "

# Row 32
$ws.Range('D32').Value = "'3.547"
$ws.Range('E32').Value = "  +2.42%  "

# Row 33
$ws.Range('D33').Value = "'3.461"
$ws.Range('E33').Value = "  +2.87%  "

# Row 34
$ws.Range('D34').Value = "'1.663"
$ws.Range('E34').Value = "  +6.49%  "

# Row 35
$ws.Range('D35').Value = "'0.9636"
$ws.Range('E35').Value = "  +1.30%  "

# Row 36
$ws.Range('D36').Value = "'2.825"
$ws.Range('E36').Value = "  +1.64%  "

# Row 37
$ws.Range('D37').Value = "'2.445"
$ws.Range('E37').Value = "  +1.78%  "

# Row 38
$ws.Range('D38').Value = "'0.5962"
$ws.Range('E38').Value = "  +4.96%  "

# Row 39
$ws.Range('D39').Value = "'0.01649"
$ws.Range('E39').Value = "  +4.20%  "

# Row 40
$ws.Range('D40').Value = "'5.929"
$ws.Range('E40').Value = "  +1.09%  "

# Row 41: Coin -> TrustWalletToken
$ws.Range('B41').Value = "TrustWalletToken"
$ws.Range('C41').Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D41').Value = "'0.8557"
$ws.Range('E41').Value = "  +2.64%  "

# Row 42: Coin -> Maker
$ws.Range('B42').Value = "Maker"
$ws.Range('C42').Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range('D42').Value = "1.057.95"
$ws.Range('E42').Value = "  +2.94%  "

# Row 43
$ws.Range('D43').Value = "'1.004"
$ws.Range('E43').Value = "  +0.06%  "

# Row 44
$ws.Range('D44').Value = "'101.41"
$ws.Range('E44').Value = "  +0.25%  "

# Row 45
$ws.Range('D45').Value = "1.868.54"
$ws.Range('E45').Value = "  +4.20%  "

# Row 46
$ws.Range('D46').Value = "0.0₈115"
$ws.Range('E46').Value = "  +20.77%  "

# Row 47
$ws.Range('D47').Value = "'59.15"

# Row 48
$ws.Range('D48').Value = "'8.225"
$ws.Range('E48').Value = "  +2.57%  "

# Row 49
$ws.Range('D49').Value = "'0.4441"
$ws.Range('E49').Value = "  +2.35%  "

# Row 50
$ws.Range('D50').Value = "'1.004"
$ws.Range('E50').Value = "  +0.39%  "

# Row 51
$ws.Range('D51').Value = "'0.05276"
$ws.Range('E51').Value = "  +1.66%  "
